$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new timesheet rows right after row 109 -----------------
# (the previously-blank separator row, and the two summary rows below it,
# all get pushed down by 2 rows to make room for the new entries)
$ws.Rows("110:111").Insert()

# New entry: 2014-04-01, 20:00-20:15
$ws.Cells.Item(110, 1).Value = 2014
$ws.Cells.Item(110, 2).Value = 4
$ws.Cells.Item(110, 3).Value = 1
$ws.Cells.Item(110, 4).Value = 0.83333333333333337
$ws.Cells.Item(110, 5).Value = 0.84375

# New entry: 2014-04-01, 21:00-22:00
$ws.Cells.Item(111, 1).Value = 2014
$ws.Cells.Item(111, 2).Value = 4
$ws.Cells.Item(111, 3).Value = 1
$ws.Cells.Item(111, 4).Value = 0.875
$ws.Cells.Item(111, 5).Value = 0.91666666666666663

# Extend the "time spent" / "time spent [h]" formulas down through the
# two newly added rows (keeps them part of the same shared formula group
# that already spans F28:F109 / G28:G109).
$ws.Range("F109:F111").Formula = "=(E109-D109)*24*60"
$ws.Range("G109:G111").Formula = "=F109/60"

# Correct the end time that was previously recorded for the 2014-04-01
# entry that is now row 109 (19:00 -> 18:45).
$ws.Cells.Item(109, 5).Value = 0.78125

# Mark H108 with the same numeric style used in column G (an otherwise
# empty helper cell used for the HELM coefficients test structure).
$ws.Cells.Item(108, 8).NumberFormat = "0.00"

# Recalculate all formulas (sums, shared formulas, etc.)
$excel.CalculateFull()

# --- Update the view so it again centers on the bottom of the sheet ----
$null = $ws.Range("A112").Select()
$excel.ActiveWindow.ScrollRow = 88
